$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency data (price, volume%, and a few reshuffled rows)
# Force text number format on target cells so numeric-looking strings
# (e.g. "69.418.42", "0.0000337", "3.00") are preserved verbatim as text,
# matching the original inline-string cell contents.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.418.42'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.70%  '
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.913.86'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.56%  '
# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.01%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '529.37'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +9.45%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.31'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.87%  '
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.614'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -1.25%  '
# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.05%  '
# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -2.91%  '
# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -4.42%  '
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000337'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -5.44%  '
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '42.13'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -3.21%  '
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.540.95'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.76%  '
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.28'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -2.21%  '
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.889.15'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.68%  '
# Row 16
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +8.37%  '
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.98'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.71%  '
# Row 18
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.63%  '
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '19.73'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.28%  '
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.364.95'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.62%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '428.05'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.66%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.38'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -5.29%  '
# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.77%  '
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '14.12'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -4.81%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.06'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +10.26%  '
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.52'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -6.34%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.62'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -4.18%  '
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '36.49'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -2.37%  '
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '13.17'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -2.67%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '672.41'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -5.50%  '
# Row 31
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -2.65%  '
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.82'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -3.09%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '68.55'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +12.00%  '
# Row 34
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = 'PEPE'
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0₃0883'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.30%  '
# Row 35
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = 'TheGraph'
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.438'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +10.49%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.93'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -2.64%  '
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '40.10'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -2.21%  '
# Row 38
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +2.13%  '
# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.06%  '
# Row 40
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.06%  '
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.23'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +4.69%  '
# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -3.77%  '
# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +7.31%  '
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.79'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -7.87%  '
# Row 45
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.22%  '
# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.95%  '
# Row 47
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'FLOKI'
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.000281'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +15.75%  '
# Row 48
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0₆0356'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +12.19%  '
# Row 49
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = 'Stacks'
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.00'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +7.06%  '
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.749.89'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +13.75%  '
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '143.79'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.36%  '
